$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 508.30768
$ws.Range("I4").Value = 759.4286
$ws.Range("J4").Value = 215.33333
$ws.Range("K4").Value = 759.4286
$ws.Range("L4").Value = 215.33333
$ws.Range("M4").Value = -645.4286
$ws.Range("N4").Value = -443.33333

$ws.Range("H10").Value = 4995
$ws.Range("J10").Value = 4995
$ws.Range("L10").Value = 4995
$ws.Range("N10").Value = -5581

$ws.Range("H44").Value = 95000
$ws.Range("J44").Value = 95000
$ws.Range("L44").Value = 95000
$ws.Range("N44").Value = -95924

$ws.Range("H74").Value = 4013.0908
$ws.Range("I74").Value = 3750
$ws.Range("J74").Value = 4232.3335
$ws.Range("K74").Value = 3750
$ws.Range("L74").Value = 4232.3335
$ws.Range("M74").Value = -2814
$ws.Range("N74").Value = -6104.3335

$ws.Range("H77").Value = 4013.0908
$ws.Range("I77").Value = 3750
$ws.Range("J77").Value = 4232.3335
$ws.Range("K77").Value = 18750
$ws.Range("L77").Value = 21161.6675
$ws.Range("M77").Value = -14070
$ws.Range("N77").Value = -30521.6675

$ws.Range("H112").Value = 6504.2334
$ws.Range("I112").Value = 940
$ws.Range("J112").Value = 7360.269
$ws.Range("K112").Value = 2820
$ws.Range("L112").Value = 22080.807
$ws.Range("M112").Value = -1712
$ws.Range("N112").Value = -24296.807

$ws.Range("H137").Value = 1029.9117
$ws.Range("I137").Value = 878.8889
$ws.Range("J137").Value = 1199.8125
$ws.Range("K137").Value = 2636.6667
$ws.Range("L137").Value = 3599.4375
$ws.Range("M137").Value = -86.66670000000022
$ws.Range("N137").Value = -8699.4375

$ws.Range("H139").Value = 67360
$ws.Range("J139").Value = 67360
$ws.Range("L139").Value = 67360
$ws.Range("N139").Value = -77640

$ws.Range("H140").Value = 85768
$ws.Range("J140").Value = 85768
$ws.Range("L140").Value = 85768
$ws.Range("N140").Value = -96128

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 1269.3636
$ws.Range("I31").Value = 1269.3636
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1269.3636
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -975.3635999999999
$ws.Range("N31").ClearContents()

$ws.Range("H61").Value = 2759.0715
$ws.Range("I61").Value = 2813.9714
$ws.Range("J61").Value = 2484.5715
$ws.Range("K61").Value = 2813.9714
$ws.Range("L61").Value = 2484.5715
$ws.Range("M61").Value = -2601.9714
$ws.Range("N61").Value = -2908.5715

$ws.Range("H74").Value = 2754.3635
$ws.Range("I74").Value = 2784.3076
$ws.Range("J74").Value = 2711.111
$ws.Range("K74").Value = 2784.3076
$ws.Range("L74").Value = 2711.111
$ws.Range("M74").Value = -1910.3076
$ws.Range("N74").Value = -4459.111

$ws.Range("H77").Value = 2754.3635
$ws.Range("I77").Value = 2784.3076
$ws.Range("J77").Value = 2711.111
$ws.Range("K77").Value = 13921.538
$ws.Range("L77").Value = 13555.555
$ws.Range("M77").Value = -9553.538
$ws.Range("N77").Value = -22291.555

$ws.Range("H132").Value = 1943.4565
$ws.Range("I132").Value = 1633.258
$ws.Range("J132").Value = 2584.5334
$ws.Range("K132").Value = 4899.774
$ws.Range("L132").Value = 7753.600199999999
$ws.Range("M132").Value = -2369.774
$ws.Range("N132").Value = -12813.6002

$ws.Range("H136").Value = 2759.0715
$ws.Range("I136").Value = 2813.9714
$ws.Range("J136").Value = 2484.5715
$ws.Range("K136").Value = 8441.914199999999
$ws.Range("L136").Value = 7453.7145
$ws.Range("M136").Value = -5891.914199999999
$ws.Range("N136").Value = -12553.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1393.238
$ws.Range("I3").Value = 669.875
$ws.Range("K3").Value = 669.875
$ws.Range("M3").Value = -555.875

$ws.Range("H138").Value = 51376
$ws.Range("J138").Value = 51376
$ws.Range("L138").Value = 51376
$ws.Range("N138").Value = -61656

$ws.Range("H140").Value = 88825
$ws.Range("J140").Value = 88825
$ws.Range("L140").Value = 88825
$ws.Range("N140").Value = -99185

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3614.5918
$ws.Range("I31").Value = 2324.7546
$ws.Range("J31").Value = 5133.7334
$ws.Range("K31").Value = 2324.7546
$ws.Range("L31").Value = 5133.7334
$ws.Range("M31").Value = -2029.7546
$ws.Range("N31").Value = -5723.7334

$ws.Range("H34").Value = 3614.5918
$ws.Range("I34").Value = 2324.7546
$ws.Range("J34").Value = 5133.7334
$ws.Range("K34").Value = 2324.7546
$ws.Range("L34").Value = 5133.7334
$ws.Range("M34").Value = -2122.7546
$ws.Range("N34").Value = -5537.7334

$ws.Range("H94").Value = 4243.85
$ws.Range("I94").Value = 800
$ws.Range("J94").Value = 4425.1055
$ws.Range("K94").Value = 800
$ws.Range("L94").Value = 4425.1055
$ws.Range("M94").Value = -349
$ws.Range("N94").Value = -5327.1055

$ws.Range("H132").Value = 1680.7446
$ws.Range("I132").Value = 1151.9143
$ws.Range("J132").Value = 3223.1667
$ws.Range("K132").Value = 3455.7429
$ws.Range("L132").Value = 9669.500100000001
$ws.Range("M132").Value = -925.7428999999997
$ws.Range("N132").Value = -14729.5001

$ws.Range("H134").Value = 1065.1936
$ws.Range("I134").Value = 899.65
$ws.Range("J134").Value = 1366.1818
$ws.Range("K134").Value = 2698.95
$ws.Range("L134").Value = 4098.5454
$ws.Range("M134").Value = -163.9499999999998
$ws.Range("N134").Value = -9168.545399999999

$ws.Range("H140").Value = 88180
$ws.Range("J140").Value = 88180
$ws.Range("L140").Value = 88180
$ws.Range("N140").Value = -98540

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H108").Value = 2800
$ws.Range("I108").Value = 2800
$ws.Range("K108").Value = 8400
$ws.Range("M108").Value = -5520

$ws.Range("H117").Value = 1842.8572
$ws.Range("I117").Value = 975
$ws.Range("K117").Value = 2925
$ws.Range("M117").Value = 517

$ws.Range("H118").Value = 2547.3572
$ws.Range("I118").Value = 1137.8
$ws.Range("K118").Value = 3413.4
$ws.Range("M118").Value = -2170.4

$ws.Range("H121").Value = 63098.47
$ws.Range("J121").Value = 69602.484
$ws.Range("L121").Value = 208807.452
$ws.Range("N121").Value = -211427.452

$ws.Range("H122").Value = 569.45
$ws.Range("I122").Value = 366
$ws.Range("J122").Value = 1179.8
$ws.Range("K122").Value = 3294
$ws.Range("L122").Value = 10618.2
$ws.Range("M122").Value = -844
$ws.Range("N122").Value = -15518.2

$ws.Range("H123").Value = 4666.6665
$ws.Range("I123").Value = 2000
$ws.Range("J123").Value = 10000
$ws.Range("K123").Value = 6000
$ws.Range("L123").Value = 30000
$ws.Range("M123").Value = -3550
$ws.Range("N123").Value = -34900

$ws.Range("H131").Value = 5155491
$ws.Range("J131").Value = 6250835.5
$ws.Range("L131").Value = 18752506.5
$ws.Range("N131").Value = -18762586.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 334433.34
$ws.Range("J40").Value = 334433.34
$ws.Range("L40").Value = 334433.34
$ws.Range("N40").Value = -334735.34

$ws.Range("H92").Value = 8091.1816
$ws.Range("J92").Value = 8091.1816
$ws.Range("L92").Value = 8091.1816
$ws.Range("N92").Value = -11835.1816

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1145
$ws.Range("I46").Value = 857.1429000000001
$ws.Range("J46").Value = 1548
$ws.Range("K46").Value = 857.1429000000001
$ws.Range("L46").Value = 1548
$ws.Range("M46").Value = -669.1429000000001
$ws.Range("N46").Value = -1924

$ws.Range("H136").Value = 2619.037
$ws.Range("I136").Value = 1739.0513
$ws.Range("J136").Value = 4907
$ws.Range("K136").Value = 5217.1539
$ws.Range("L136").Value = 14721
$ws.Range("M136").Value = -2667.1539
$ws.Range("N136").Value = -19821

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1435.9286
$ws.Range("I126").Value = 1040.5555
$ws.Range("J126").Value = 2147.6
$ws.Range("K126").Value = 3121.6665
$ws.Range("L126").Value = 6442.799999999999
$ws.Range("M126").Value = -651.6664999999998
$ws.Range("N126").Value = -11382.8

$ws.Range("H136").Value = 792.5854
$ws.Range("I136").Value = 662.62067
$ws.Range("J136").Value = 1106.6666
$ws.Range("K136").Value = 1987.86201
$ws.Range("L136").Value = 3319.9998
$ws.Range("M136").Value = 562.1379899999999
$ws.Range("N136").Value = -8419.9998
